$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-27
$values = @(
    @(6, 7),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(1, 3),
    @(2, 5),
    @(8, 9),
    @(7, 7),
    @(8, 9),
    @(10, 10),
    @(4, 4),
    @(6, 7),
    @(5, 6),
    @(6, 8),
    @(7, 8),
    @(8, 9),
    @(6, 6),
    @(5, 6),
    @(9, 9),
    @(5, 5),
    @(7, 7),
    @(6, 6),
    @(2, 2),
    @(6, 6),
    @(2, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
